$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Target data for rows 2-13 (columns A-T), reflecting the NATMI re-run after
# adding the "M2" cell type / cluster to the analysis (per Dr Hou's advice).
# Column order: A=Sending cluster, B=Ligand symbol, C=Receptor symbol,
# D=Target cluster, E..T = numeric stats.
$data = @(
    @('ECs','Efna1','Epha4','ECs',3.0,1.0,23.60223533333333,70.80670599999999,0.8824726436021215,0.8824726436021214,3.0,1.0,7.839997333333334,23.519992,0.3930770090432645,0.3930770090432645,185.0414620740391,1665.373158666352,0.3468797073096246,0.3468797073096246),
    @('ECs','Efna1','Epha4','FAPs',3.0,1.0,23.60223533333333,70.80670599999999,0.8824726436021215,0.8824726436021214,3.0,1.0,11.57455166666667,34.723655,0.5803178185796234,0.5803178185796235,273.1852923144922,2458.66763083043,0.5121145994913766,0.5121145994913766),
    @('ECs','Efna1','Epha4','M2',3.0,1.0,23.60223533333333,70.80670599999999,0.8824726436021215,0.8824726436021214,2.0,0.6666666666666666,0.04482366666666667,0.134471,0.002247341686300608,0.002247341686300608,1.057938729169555,9.521448562525999,0.001983217558986948,0.001983217558986947),
    @('ECs','Efna1','Epha4','sCs',3.0,1.0,23.60223533333333,70.80670599999999,0.8824726436021215,0.8824726436021214,3.0,1.0,0.4858216666666667,1.457465,0.02435783069081152,0.02435783069081152,11.46647730669889,103.19829576029,0.02149511924213333,0.02149511924213333),
    @('FAPs','Efna1','Epha4','ECs',3.0,1.0,2.633202,7.899606,0.09845375648791208,0.09845375648791205,3.0,1.0,7.839997333333334,23.519992,0.3930770090432645,0.3930770090432645,20.644296658128,185.798669923152,0.03869990812934238,0.03869990812934236),
    @('FAPs','Efna1','Epha4','FAPs',3.0,1.0,2.633202,7.899606,0.09845375648791208,0.09845375648791205,3.0,1.0,11.57455166666667,34.723655,0.5803178185796234,0.5803178185796235,30.47813259777,274.30319337993,0.05713446919603458,0.05713446919603457),
    @('FAPs','Efna1','Epha4','M2',3.0,1.0,2.633202,7.899606,0.09845375648791208,0.09845375648791205,2.0,0.6666666666666666,0.04482366666666667,0.134471,0.002247341686300608,0.002247341686300608,0.118029768714,1.062267918426,0.0002212592311281738,0.0002212592311281737),
    @('FAPs','Efna1','Epha4','sCs',3.0,1.0,2.633202,7.899606,0.09845375648791208,0.09845375648791205,3.0,1.0,0.4858216666666667,1.457465,0.02435783069081152,0.02435783069081152,1.27926658431,11.51339925879,0.002398119931406948,0.002398119931406948),
    @('sCs','Efna1','Epha4','ECs',3.0,1.0,0.5101343333333334,1.530403,0.01907359990996641,0.0190735999099664,3.0,1.0,7.839997333333334,23.519992,0.3930770090432645,0.3930770090432645,3.999451812975111,35.995066316776,0.007497393604297475,0.007497393604297471),
    @('sCs','Efna1','Epha4','FAPs',3.0,1.0,0.5101343333333334,1.530403,0.01907359990996641,0.0190735999099664,3.0,1.0,11.57455166666667,34.723655,0.5803178185796234,0.5803178185796235,5.904576198107223,53.141185782965,0.01106874989221221,0.0110687498922122),
    @('sCs','Efna1','Epha4','M2',3.0,1.0,0.5101343333333334,1.530403,0.01907359990996641,0.0190735999099664,2.0,0.6666666666666666,0.04482366666666667,0.134471,0.002247341686300608,0.002247341686300608,0.02286609131255556,0.205794821813,0.00004286489618548703,0.00004286489618548702),
    @('sCs','Efna1','Epha4','sCs',3.0,1.0,0.5101343333333334,1.530403,0.01907359990996641,0.0190735999099664,3.0,1.0,0.4858216666666667,1.457465,0.02435783069081152,0.02435783069081152,0.2478343120438889,2.230508808395,0.0004645915172712396,0.0004645915172712394)
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $data[$i]
    $r = $i + 2
    for ($c = 0; $c -lt $row.Length; $c++) {
        $ws.Cells.Item($r, $c + 1).Value = $row[$c]
    }
}
